$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new date-like labels to be stored as plain text (not auto-converted
# to Excel date serials), matching how the existing "Serie" column is stored.
$ws.Range("A173:A174").NumberFormat = "@"

$ws.Range("A173").Value = "07-09-2021"
$ws.Range("B173").Value = -0.29
$ws.Range("C173").Value = -0.16
$ws.Range("D173").Value = 0.01

$ws.Range("A174").Value = "08-09-2021"
$ws.Range("B174").Value = -0.24
$ws.Range("C174").Value = -0.07000000000000001
$ws.Range("D174").Value = 0.07000000000000001

# Restore the default (unstyled) cell formatting so the new rows match the
# look of the rest of the table.
$ws.Range("A173:A174").Style = "Normal"
